$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "EX_glc_LPAREN_e_RPAREN_",
    "EX_h2o_LPAREN_e_RPAREN_",
    "EX_h_LPAREN_e_RPAREN_",
    "EX_inost_LPAREN_e_RPAREN_",
    "EX_k_LPAREN_e_RPAREN_",
    "EX_na1_LPAREN_e_RPAREN_",
    "EX_nh4_LPAREN_e_RPAREN_",
    "EX_o2_LPAREN_e_RPAREN_",
    "EX_pi_LPAREN_e_RPAREN_",
    "EX_so4_LPAREN_e_RPAREN_",
    "trehalose_c_tp"
)

$values = @(
    0.08000000000000002,
    41.65422506566294,
    41.65422506566294,
    41.64483387653971,
    41.65422506566294,
    41.65422506566294,
    19.45251739130443,
    0.08000000000000002,
    19.45251739130411,
    19.45251739130492,
    41.65422506566343
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
